$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1262.7391
$ws.Range("I70").Value = 653.3333
$ws.Range("J70").Value = 1354.15
$ws.Range("K70").Value = 1959.9999
$ws.Range("L70").Value = 4062.45
$ws.Range("M70").Value = -1689.9999
$ws.Range("N70").Value = -4602.450000000001
$ws.Range("H73").Value = 1262.7391
$ws.Range("I73").Value = 653.3333
$ws.Range("J73").Value = 1354.15
$ws.Range("K73").Value = 1959.9999
$ws.Range("L73").Value = 4062.45
$ws.Range("M73").Value = -1023.9999
$ws.Range("N73").Value = -5934.450000000001
$ws.Range("H76").Value = 10924.167
$ws.Range("I76").Value = 2737.5
$ws.Range("J76").Value = 27297.5
$ws.Range("K76").Value = 2737.5
$ws.Range("L76").Value = 27297.5
$ws.Range("M76").Value = -2422.5
$ws.Range("N76").Value = -27927.5
$ws.Range("H79").Value = 10924.167
$ws.Range("I79").Value = 2737.5
$ws.Range("J79").Value = 27297.5
$ws.Range("K79").Value = 2737.5
$ws.Range("L79").Value = 27297.5
$ws.Range("M79").Value = -1645.5
$ws.Range("N79").Value = -29481.5
$ws.Range("H127").Value = 1755.0869
$ws.Range("I127").Value = 1360.25
$ws.Range("K127").Value = 4080.75
$ws.Range("M127").Value = 879.25
$ws.Range("H141").Value = 3551.1765
$ws.Range("I141").Value = 1554.7826
$ws.Range("J141").Value = 7725.4546
$ws.Range("K141").Value = 4664.3478
$ws.Range("L141").Value = 23176.3638
$ws.Range("M141").Value = 515.6522000000004
$ws.Range("N141").Value = -33536.3638

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1373.5172
$ws.Range("I45").Value = 1217.1428
$ws.Range("K45").Value = 1217.1428
$ws.Range("M45").Value = -840.1428000000001
$ws.Range("H61").Value = 1463.4849
$ws.Range("I61").Value = 739.5
$ws.Range("J61").Value = 2332.2666
$ws.Range("K61").Value = 739.5
$ws.Range("L61").Value = 2332.2666
$ws.Range("M61").Value = -527.5
$ws.Range("N61").Value = -2756.2666
$ws.Range("H63").Value = 2267.5325
$ws.Range("I63").Value = 2258.5715
$ws.Range("J63").Value = 2357.1428
$ws.Range("K63").Value = 2258.5715
$ws.Range("L63").Value = 2357.1428
$ws.Range("M63").Value = -1572.5715
$ws.Range("N63").Value = -3729.1428
$ws.Range("H66").Value = 2267.5325
$ws.Range("I66").Value = 2258.5715
$ws.Range("J66").Value = 2357.1428
$ws.Range("K66").Value = 11292.8575
$ws.Range("L66").Value = 11785.714
$ws.Range("M66").Value = -7860.8575
$ws.Range("N66").Value = -18649.714
$ws.Range("H74").Value = 1130.0714
$ws.Range("I74").Value = 613.73334
$ws.Range("J74").Value = 1725.8462
$ws.Range("K74").Value = 613.73334
$ws.Range("L74").Value = 1725.8462
$ws.Range("M74").Value = 260.26666
$ws.Range("N74").Value = -3473.8462
$ws.Range("H77").Value = 1130.0714
$ws.Range("I77").Value = 613.73334
$ws.Range("J77").Value = 1725.8462
$ws.Range("K77").Value = 3068.6667
$ws.Range("L77").Value = 8629.231
$ws.Range("M77").Value = 1299.3333
$ws.Range("N77").Value = -17365.231
$ws.Range("H88").Value = 2666.6667
$ws.Range("I88").Value = 2666.6667
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 2666.6667
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = -2260.6667
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 2666.6667
$ws.Range("I91").Value = 2666.6667
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 2666.6667
$ws.Range("L91").Value = 0
$ws.Range("M91").Value = -1262.6667
$ws.Range("N91").ClearContents()
$ws.Range("H132").Value = 1460.28
$ws.Range("I132").Value = 529.8570999999999
$ws.Range("J132").Value = 3631.2666
$ws.Range("K132").Value = 1589.5713
$ws.Range("L132").Value = 10893.7998
$ws.Range("M132").Value = 940.4287000000002
$ws.Range("N132").Value = -15953.7998
$ws.Range("H136").Value = 1463.4849
$ws.Range("I136").Value = 739.5
$ws.Range("J136").Value = 2332.2666
$ws.Range("K136").Value = 2218.5
$ws.Range("L136").Value = 6996.7998
$ws.Range("M136").Value = 331.5
$ws.Range("N136").Value = -12096.7998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1939.7778
$ws.Range("I86").Value = 1994
$ws.Range("J86").Value = 1750
$ws.Range("K86").Value = 1994
$ws.Range("L86").Value = 1750
$ws.Range("M86").Value = -871
$ws.Range("N86").Value = -3996
$ws.Range("H89").Value = 1939.7778
$ws.Range("I89").Value = 1994
$ws.Range("J89").Value = 1750
$ws.Range("K89").Value = 9970
$ws.Range("L89").Value = 8750
$ws.Range("M89").Value = -4354
$ws.Range("N89").Value = -19982
$ws.Range("H134").Value = 1560.0256
$ws.Range("I134").Value = 1411.9706
$ws.Range("J134").Value = 2566.8
$ws.Range("K134").Value = 4235.9118
$ws.Range("L134").Value = 7700.400000000001
$ws.Range("M134").Value = -1700.9118
$ws.Range("N134").Value = -12770.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 9486.571
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 9486.571
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 9486.571
$ws.Range("N11").Value = -9766.571
$ws.Range("M11").ClearContents()
$ws.Range("H25").Value = 10100
$ws.Range("I25").Value = 900
$ws.Range("K25").Value = 900
$ws.Range("M25").Value = -726
$ws.Range("H105").Value = 1142.05
$ws.Range("I105").Value = 940.7692
$ws.Range("J105").Value = 1515.8572
$ws.Range("K105").Value = 940.7692
$ws.Range("L105").Value = 1515.8572
$ws.Range("M105").Value = 806.2308
$ws.Range("N105").Value = -5009.8572

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1126.076
$ws.Range("I68").Value = 741.28
$ws.Range("J68").Value = 1584.1666
$ws.Range("K68").Value = 2223.84
$ws.Range("L68").Value = 4752.4998
$ws.Range("M68").Value = -1412.84
$ws.Range("N68").Value = -6374.4998
$ws.Range("H71").Value = 1126.076
$ws.Range("I71").Value = 741.28
$ws.Range("J71").Value = 1584.1666
$ws.Range("K71").Value = 6671.52
$ws.Range("L71").Value = 14257.4994
$ws.Range("M71").Value = -2615.52
$ws.Range("N71").Value = -22369.4994
$ws.Range("H131").Value = 12500937
$ws.Range("I131").Value = 55556068
$ws.Range("J131").Value = 1060.2903
$ws.Range("K131").Value = 166668204
$ws.Range("L131").Value = 3180.8709
$ws.Range("M131").Value = -166663164
$ws.Range("N131").Value = -13260.8709

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5433.222
$ws.Range("I70").Value = 4908.25
$ws.Range("J70").Value = 6483.1665
$ws.Range("K70").Value = 4908.25
$ws.Range("L70").Value = 6483.1665
$ws.Range("M70").Value = -4638.25
$ws.Range("N70").Value = -7023.1665
$ws.Range("H73").Value = 5433.222
$ws.Range("I73").Value = 4908.25
$ws.Range("J73").Value = 6483.1665
$ws.Range("K73").Value = 4908.25
$ws.Range("L73").Value = 6483.1665
$ws.Range("M73").Value = -3972.25
$ws.Range("N73").Value = -8355.166499999999
$ws.Range("H113").Value = 4593.706
$ws.Range("I113").Value = 5814.5654
$ws.Range("J113").Value = 2041
$ws.Range("K113").Value = 5814.5654
$ws.Range("L113").Value = 2041
$ws.Range("M113").Value = -3644.5654
$ws.Range("N113").Value = -6381

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2538.6924
$ws.Range("I132").Value = 1442.409
$ws.Range("J132").Value = 3957.4119
$ws.Range("K132").Value = 4327.227000000001
$ws.Range("L132").Value = 11872.2357
$ws.Range("M132").Value = -1797.227000000001
$ws.Range("N132").Value = -16932.2357
$ws.Range("H136").Value = 2416.0322
$ws.Range("I136").Value = 1220.6
$ws.Range("J136").Value = 2985.2856
$ws.Range("K136").Value = 3661.8
$ws.Range("L136").Value = 8955.856800000001
$ws.Range("M136").Value = -1111.8
$ws.Range("N136").Value = -14055.8568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1588.5294
$ws.Range("I96").Value = 1422.8889
$ws.Range("J96").Value = 1774.875
$ws.Range("K96").Value = 1422.8889
$ws.Range("L96").Value = 1774.875
$ws.Range("M96").Value = -49.88889999999992
$ws.Range("N96").Value = -4520.875
$ws.Range("H107").Value = 1288.8889
$ws.Range("I107").Value = 1400
$ws.Range("J107").Value = 1150
$ws.Range("K107").Value = 4200
$ws.Range("L107").Value = 3450
$ws.Range("M107").Value = -2280
$ws.Range("N107").Value = -7290
$ws.Range("H136").Value = 8468.3125
$ws.Range("I136").Value = 2054.2222
$ws.Range("J136").Value = 16715
$ws.Range("K136").Value = 6162.6666
$ws.Range("L136").Value = 50145
$ws.Range("M136").Value = -3612.6666
$ws.Range("N136").Value = -55245
